$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 83.46154
$ws.Range("I8").Value = 83.46154
$ws.Range("K8").Value = 250.38462
$ws.Range("M8").Value = -111.38462
$ws.Range("H15").Value = 1928.9753
$ws.Range("I15").Value = 1928.9753
$ws.Range("K15").Value = 5786.9259
$ws.Range("M15").Value = -5617.9259
$ws.Range("H76").Value = 11500281
$ws.Range("I76").Value = 7395.7144
$ws.Range("J76").Value = 41669104
$ws.Range("K76").Value = 7395.7144
$ws.Range("L76").Value = 41669104
$ws.Range("M76").Value = -7080.7144
$ws.Range("N76").Value = -41669734
$ws.Range("H79").Value = 11500281
$ws.Range("I79").Value = 7395.7144
$ws.Range("J79").Value = 41669104
$ws.Range("K79").Value = 7395.7144
$ws.Range("L79").Value = 41669104
$ws.Range("M79").Value = -6303.7144
$ws.Range("N79").Value = -41671288
$ws.Range("H111").Value = 500
$ws.Range("I111").Value = 500
$ws.Range("K111").Value = 1500
$ws.Range("M111").Value = 1567
$ws.Range("H116").Value = 2480
$ws.Range("J116").Value = 2480
$ws.Range("L116").Value = 2480
$ws.Range("N116").Value = -9364
$ws.Range("H135").Value = 910.13794
$ws.Range("I135").Value = 702.2041
$ws.Range("J135").Value = 2042.2222
$ws.Range("K135").Value = 6319.8369
$ws.Range("L135").Value = 18379.9998
$ws.Range("M135").Value = -3784.8369
$ws.Range("N135").Value = -23449.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2371.3572
$ws.Range("J3").Value = 2530.7693
$ws.Range("L3").Value = 2530.7693
$ws.Range("N3").Value = -2760.7693
$ws.Range("H132").Value = 19609480
$ws.Range("I132").Value = 27779274
$ws.Range("J132").Value = 3269890
$ws.Range("K132").Value = 83337822
$ws.Range("L132").Value = 9809670
$ws.Range("M132").Value = -83335292
$ws.Range("N132").Value = -9814730

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2586660.2
$ws.Range("I86").Value = 3367
$ws.Range("J86").Value = 5815777
$ws.Range("K86").Value = 3367
$ws.Range("L86").Value = 5815777
$ws.Range("M86").Value = -2244
$ws.Range("N86").Value = -5818023
$ws.Range("H89").Value = 2586660.2
$ws.Range("I89").Value = 3367
$ws.Range("J89").Value = 5815777
$ws.Range("K89").Value = 16835
$ws.Range("L89").Value = 29078885
$ws.Range("M89").Value = -11219
$ws.Range("N89").Value = -29090117

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1346.3914
$ws.Range("I31").Value = 1006.9545
$ws.Range("J31").Value = 1657.5416
$ws.Range("K31").Value = 1006.9545
$ws.Range("L31").Value = 1657.5416
$ws.Range("M31").Value = -711.9545
$ws.Range("N31").Value = -2247.5416
$ws.Range("H34").Value = 1346.3914
$ws.Range("I34").Value = 1006.9545
$ws.Range("J34").Value = 1657.5416
$ws.Range("K34").Value = 1006.9545
$ws.Range("L34").Value = 1657.5416
$ws.Range("M34").Value = -804.9545
$ws.Range("N34").Value = -2061.5416
$ws.Range("H58").Value = 13514355
$ws.Range("I58").Value = 19231600
$ws.Range("J58").Value = 865.3182
$ws.Range("K58").Value = 19231600
$ws.Range("L58").Value = 865.3182
$ws.Range("M58").Value = -19231397
$ws.Range("N58").Value = -1271.3182
$ws.Range("H136").Value = 13514355
$ws.Range("I136").Value = 19231600
$ws.Range("J136").Value = 865.3182
$ws.Range("K136").Value = 57694800
$ws.Range("L136").Value = 2595.9546
$ws.Range("M136").Value = -57692250
$ws.Range("N136").Value = -7695.9546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4666.6924
$ws.Range("I3").Value = 4583.375
$ws.Range("J3").Value = 4800
$ws.Range("K3").Value = 13750.125
$ws.Range("L3").Value = 14400
$ws.Range("M3").Value = -13638.125
$ws.Range("N3").Value = -14624
$ws.Range("H10").Value = 106.181816
$ws.Range("I10").Value = 106.181816
$ws.Range("K10").Value = 318.545448
$ws.Range("M10").Value = -179.545448
$ws.Range("H23").Value = 208.5
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 244.66667
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 734.00001
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -1204.00001
$ws.Range("H68").Value = 6331.353
$ws.Range("I68").Value = 380
$ws.Range("J68").Value = 8811.083
$ws.Range("K68").Value = 1140
$ws.Range("L68").Value = 26433.249
$ws.Range("M68").Value = -329
$ws.Range("N68").Value = -28055.249
$ws.Range("H71").Value = 6331.353
$ws.Range("I71").Value = 380
$ws.Range("J71").Value = 8811.083
$ws.Range("K71").Value = 3420
$ws.Range("L71").Value = 79299.747
$ws.Range("M71").Value = 636
$ws.Range("N71").Value = -87411.747
$ws.Range("H122").Value = 10969421
$ws.Range("I122").Value = 62500264
$ws.Range("J122").Value = 5412.447
$ws.Range("K122").Value = 562502376
$ws.Range("L122").Value = 48712.023
$ws.Range("M122").Value = -562499926
$ws.Range("N122").Value = -53612.023
$ws.Range("H131").Value = 837.5263
$ws.Range("I131").Value = 446.27274
$ws.Range("J131").Value = 888.7619
$ws.Range("K131").Value = 1338.81822
$ws.Range("L131").Value = 2666.2857
$ws.Range("M131").Value = 3701.18178
$ws.Range("N131").Value = -12746.2857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1187.4
$ws.Range("I97").Value = 962.5
$ws.Range("J97").Value = 1524.75
$ws.Range("K97").Value = 962.5
$ws.Range("L97").Value = 1524.75
$ws.Range("M97").Value = -466.5
$ws.Range("N97").Value = -2516.75
$ws.Range("H132").Value = 26519.111
$ws.Range("I132").Value = 18943
$ws.Range("K132").Value = 56829
$ws.Range("M132").Value = -54299

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1983.7273
$ws.Range("I22").Value = 500.5
$ws.Range("J22").Value = 2313.3333
$ws.Range("K22").Value = 500.5
$ws.Range("L22").Value = 2313.3333
$ws.Range("M22").Value = -205.5
$ws.Range("N22").Value = -2903.3333
$ws.Range("H27").Value = 1983.7273
$ws.Range("I27").Value = 500.5
$ws.Range("J27").Value = 2313.3333
$ws.Range("K27").Value = 500.5
$ws.Range("L27").Value = 2313.3333
$ws.Range("M27").Value = -393.5
$ws.Range("N27").Value = -2527.3333
$ws.Range("H46").Value = 1564.381
$ws.Range("I46").Value = 551.6667
$ws.Range("J46").Value = 2323.9167
$ws.Range("K46").Value = 551.6667
$ws.Range("L46").Value = 2323.9167
$ws.Range("M46").Value = -363.6667
$ws.Range("N46").Value = -2699.9167
$ws.Range("H132").Value = 27783254
$ws.Range("I132").Value = 47621628
$ws.Range("J132").Value = 9530.667
$ws.Range("K132").Value = 142864884
$ws.Range("L132").Value = 28592.001
$ws.Range("M132").Value = -142862354
$ws.Range("N132").Value = -33652.001
$ws.Range("H136").Value = 25453106
$ws.Range("I136").Value = 11075718
$ws.Range("J136").Value = 66668290
$ws.Range("K136").Value = 33227154
$ws.Range("L136").Value = 200004870
$ws.Range("M136").Value = -33224604
$ws.Range("N136").Value = -200009970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26187.957
$ws.Range("I132").Value = 37131.207
$ws.Range("K132").Value = 111393.621
$ws.Range("M132").Value = -108863.621
$ws.Range("H136").Value = 9438503
$ws.Range("I136").Value = 13894559
$ws.Range("J136").Value = 2150.2942
$ws.Range("K136").Value = 41683677
$ws.Range("L136").Value = 6450.882599999999
$ws.Range("M136").Value = -41681127
$ws.Range("N136").Value = -11550.8826
